$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9016587138175964
$ws.Range("B1").Value = 1.403303265571594
$ws.Range("C1").Value = 4.817044258117676
$ws.Range("D1").Value = 3.018893241882324
$ws.Range("E1").Value = 0.4381093680858612
